$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "CodeModule" column used to be the last column (F); it is now moved to
# the front (A), shifting ModuleName..Coordinator one column to the right,
# and the module codes are renumbered from the AP2x series to the AP3x series.
$data = @(
    @("CodeModule", "ModuleName", "ElementName1", "ElementName2", "Dept_Attachement", "Coordinator"),
    @("AP31", "MATHS3", "ANALYSE3", "ALGEBRE3", "SIC", "EZZINE"),
    @("AP32", "PHYSIQUE4", "Electromagnetisme", "Electo", "GEI", "chahboun"),
    @("AP33", "MECANIQUE2", "Mecanique industrielle", "mecanique_ind", "GEI", "SARSRI"),
    @("AP34", "CHIMIE", "Atomistique", "Chimie Organique", "SIC", "TOUHAMI"),
    @("AP35", "LC3", "Francais1", "Activites d'ouverture", "SIC", "HARIS"),
    @("AP36", "MATHS4", "Algebre4", "Analyse4", "SIC", "AMAL"),
    @("AP37", "PYSIQUE5", "Electronique", "instrumentation2", "SIC", "Britel"),
    @("AP38", "PYSIQUE5", "Thermodynamique", "instrumentation3", "SIC", "FILALI"),
    @("AP39", "INFO1", "Programmation C", "Algorithme", "SIC", "ALAMI")
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r + 1, $c + 1).Value = $row[$c]
    }
}

# Update the active selection as recorded in the saved workbook view.
$ws.Range("E7").Select()
